# "created clean all function" — fill in the previously-blank Sales
# Amount / Value cells for Cheddar Bacon and Crispy Bacon Burger, refresh
# the Double Bacon row, and recompute the Total row so it cleanly sums
# the four product rows above it.
#
# Target table:
#   Classic Burger        3     38.70
#   Cheddar Bacon          6     95.40
#   Double Bacon          13    245.70
#   Crispy Bacon Burger    5     89.50
#   Total                 27    469.30
#
# These numbers must land as text (shared-string) cells, matching the
# existing "3" / "38.70" style cells already on the sheet, rather than
# as native numeric cells (which is what a plain Range.Value = "6"
# assignment would auto-coerce to, and which would also force a new
# cell style). Writing a text formula and then pasting-special just the
# values back on top of itself converts the formula result into a
# literal string cell while keeping the original style untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextNumber($addr, $text) {
    $ws.Range($addr).Formula = "=""" + $text + """"
    $ws.Range($addr).Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-TextNumber "B3" "6"
Set-TextNumber "C3" "95.40"

Set-TextNumber "B4" "13"
Set-TextNumber "C4" "245.70"

Set-TextNumber "B5" "5"
Set-TextNumber "C5" "89.50"

Set-TextNumber "B6" "27"
Set-TextNumber "C6" "469.30"
